# #5: insurance, claim, debt, investment done
#
# The "債務" (debt) worksheet (sheet7) is expanded from 7 columns (A:G) to
# 14 columns (A:N). Row 1 switches from holding (duplicated) data values to
# proper field-name headers, and row 2 keeps/gets the real record values,
# now including the extra property/legislator metadata columns that every
# other sheet in this workbook already has
# (property_category/category/date/legislator_name/legislator_id/
#  source_file/index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# ---- Row 1: column headers (B1:N1) ----
$ws.Cells.Item(1, 2).Value  = "species"
$ws.Cells.Item(1, 3).Value  = "debtor"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "total"
$ws.Cells.Item(1, 6).Value  = "register_date"
$ws.Cells.Item(1, 7).Value  = "register_reason"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Give the new H1:N1 header cells the same bold / centred / bordered look
# already used by B1:G1 (cellXfs style index 1 in the original workbook) by
# cloning the format from B1 in one shot, instead of rebuilding it property
# by property (which would leave half-built styles behind in styles.xml).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---- Row 2: the actual debt record (A2:N2) ----
# A2 (125), B2:G2 (房屋貸款/王怡心/國泰世華.../12029020/96年06月20日/買房子)
# already hold the right values - only the new H2:N2 columns need filling.
$ws.Cells.Item(2, 8).Value  = "debt"
$ws.Cells.Item(2, 9).Value  = "normal"

# J2 ("date") is the literal text "2012-04-20"; pre-format as Text so Excel
# doesn't reinterpret the ISO-looking string as a date serial.
$dateCell = $ws.Cells.Item(2, 10)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2012-04-20"

$ws.Cells.Item(2, 11).Value = "費鴻泰"
$ws.Cells.Item(2, 12).Value = 1365
$ws.Cells.Item(2, 13).Value = "tmpe52e1"
$ws.Cells.Item(2, 14).Value = 125

$excel.CutCopyMode = $false
